$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-6 (columns D, L, M, N, O, P, R, S)
# This is effectively a cyclic rotation of the weekly records:
# old row5 -> row2, old row6 -> row3, old row2 -> row4, old row3 -> row5, old row4 -> row6

$rows = @(
    @{ Row = 2; D = 45072; L = "Segunda"; M = 100; N = 16000; O = 16000; P = 16000; R = "Provincia de Chacabuco"; S = 889 },
    @{ Row = 3; D = 45072; L = "Segunda"; M = 100; N = 17000; O = 17000; P = 17000; R = "Provincia de Limarí";   S = 944 },
    @{ Row = 4; D = 44252; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; R = "Región Metropolitana";  S = 750 },
    @{ Row = 5; D = 44250; L = "Primera"; M = 200; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";  S = 806 },
    @{ Row = 6; D = 44253; L = "Primera"; M = 160; N = 14000; O = 15000; P = 14500; R = "Región Metropolitana";  S = 806 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($rowNum, 12).Value = $r.L   # L: Calidad
    $ws.Cells.Item($rowNum, 13).Value = $r.M   # M: Volumen
    $ws.Cells.Item($rowNum, 14).Value = $r.N   # N: Precio mínimo
    $ws.Cells.Item($rowNum, 15).Value = $r.O   # O: Precio máximo
    $ws.Cells.Item($rowNum, 16).Value = $r.P   # P: Precio promedio ponderado
    $ws.Cells.Item($rowNum, 18).Value = $r.R   # R: Origen
    $ws.Cells.Item($rowNum, 19).Value = $r.S   # S: Precio $/Kg
}
